$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column H ("Absent") should be the complement of column E ("Real"):
# a student is Absent on a date if they were not marked Real/present.
$lastRow = 21
for ($r = 3; $r -le $lastRow; $r++) {
    $real = $ws.Cells.Item($r, 5).Value()  # column E = Real
    if ($real -eq 1) {
        $ws.Cells.Item($r, 8).Value = 0
    } else {
        $ws.Cells.Item($r, 8).Value = 1
    }
}
